# "story and then some" -- mark the last three list items of the outline
# ("class utility", "line shape (2.001)", "curve shape (2)") as done, the
# same way earlier items in this document were marked done: strike through
# the whole paragraph (paragraph mark included, so the strike "sticks" to
# the bullet/line itself, matching List Paragraph items elsewhere in the
# doc that already carry <w:rPr><w:strike/></w:rPr> on both the pPr and
# every run).
#
# The "line shape (2.001)" paragraph also carries the document's hidden
# _GoBack bookmark stuck in the middle of the number ("2.00" + bookmark +
# "1)"). In the target revision that bookmark has been moved out of the
# middle of the text to the very end of the paragraph (after "2.001)"),
# so the number reads as a single uninterrupted run sequence again before
# the strike-through is applied.

$d = $word.ActiveDocument

function Get-ParaByText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs.Item($i)
        if ($cand.Range.Text.TrimEnd("`r", "`a") -eq $text) {
            return $i
        }
    }
    throw "paragraph not found: $text"
}

$classUtilityIdx = Get-ParaByText("class utility")
$lineShapeIdx    = Get-ParaByText("line shape (2.001)")
$curveShapeIdx   = Get-ParaByText("curve shape (2)")

# --- "class utility" : just strike the whole paragraph -------------------
$d.Paragraphs.Item($classUtilityIdx).Range.Font.StrikeThrough = 1

# --- "line shape (2.001)" : relocate the _GoBack bookmark to the end -----
# first, then strike the whole paragraph.
$goBack = $d.Bookmarks.Item("_GoBack")
$p = $d.Paragraphs.Item($lineShapeIdx)

# Cut the text that currently sits *after* the bookmark ("1)") ...
$tailRange = $d.Range($goBack.End, $p.Range.End - 1)
$tailRange.Cut()

# ... drop the now-empty bookmark ...
$d.Bookmarks.Item("_GoBack").Delete()

# ... and paste "1)" back immediately, restoring "2.001)" with no bookmark
# in the middle any more.
$pasteAt = $d.Paragraphs.Item($lineShapeIdx).Range.End - 1
$d.Range($pasteAt, $pasteAt).Paste()

# Appending a bookmark exactly at "paragraph end minus one" (the very last
# character slot before the pilcrow) isn't reliable here, so park a throw-
# away character there first, anchor the bookmark just *before* it (a safe,
# interior position), and then delete the placeholder again - leaving the
# bookmark as the last thing in the paragraph, right after ")".
$p = $d.Paragraphs.Item($lineShapeIdx)
$endPos = $p.Range.End - 1
$d.Range($endPos, $endPos).InsertAfter("X")

$p = $d.Paragraphs.Item($lineShapeIdx)
$bookmarkPos = $p.Range.End - 2
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))

$p = $d.Paragraphs.Item($lineShapeIdx)
$placeholderPos = $p.Range.End - 1
$d.Range($placeholderPos - 1, $placeholderPos).Delete()

$d.Paragraphs.Item($lineShapeIdx).Range.Font.StrikeThrough = 1

# --- "curve shape (2)" : just strike the whole paragraph -----------------
$d.Paragraphs.Item($curveShapeIdx).Range.Font.StrikeThrough = 1
